# Update Reaction_number values (column C, rows 2-20) on both the NBR and BAR sheets.

$wb = $excel.ActiveWorkbook

$wsNBR = $wb.Worksheets.Item("NBR")
$wsBAR = $wb.Worksheets.Item("BAR")

$nbrValues = @(524, 529, 520, 515, 507, 498, 500, 485, 485, 483, 461, 456, 454, 453, 450, 441, 439, 435, 434)
$barValues = @(659, 638, 637, 639, 642, 631, 628, 630, 631, 625, 628, 635, 632, 627, 630, 629, 632, 634, 631)

for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNBR.Cells.Item($row, 3).Value = $nbrValues[$i]
}

for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $wsBAR.Cells.Item($row, 3).Value = $barValues[$i]
}
